# Actualizacion - Semana 8
# Update the "Produccion" column (G) values with the newly computed
# weekly production figures. Values are grouped in blocks of rows that
# previously shared the same figure; each block is updated to its new
# value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Rows = 2..6;   Value = 4682.314840316323 },
    @{ Rows = 7..10;  Value = 5485.172645633665 },
    @{ Rows = 11..14; Value = 4716.307814546833 },
    @{ Rows = 15..18; Value = 6595.819040059065 },
    @{ Rows = 19..21; Value = 4999.063386440928 },
    @{ Rows = 22..25; Value = 7447.960279713438 },
    @{ Rows = 26..30; Value = 5444.06118687639 },
    @{ Rows = 31..34; Value = 4446.212064655101 },
    @{ Rows = 35..39; Value = 3504.024688385195 },
    @{ Rows = 40..43; Value = 5385.679067766575 },
    @{ Rows = 44..47; Value = 7388.869513622338 },
    @{ Rows = 48..51; Value = 7582.515471984145 },
    @{ Rows = 52..55; Value = 9514.654555849344 },
    @{ Rows = 56..59; Value = 5410.030976527181 },
    @{ Rows = 60..64; Value = 4485.554779177548 },
    @{ Rows = 65..68; Value = 3731.890260596616 },
    @{ Rows = 69..72; Value = 6421.037168220616 },
    @{ Rows = 73..77; Value = 6389.296100741921 },
    @{ Rows = 78..80; Value = 4222.840907883964 },
    @{ Rows = 81..82; Value = 3199.757991717083 }
)

foreach ($u in $updates) {
    foreach ($r in $u.Rows) {
        $ws.Cells.Item($r, 7).Value = $u.Value
    }
}
